# Fruta / hortaliza, semanal
# A new weekly price record for "Membrillo" (Vega Modelo de Temuco) is added
# right after the most recent existing record (row 178), pushing all
# subsequent historical rows down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 179; this shifts rows 179:247 down to 180:248
# and keeps their formatting/styles intact.
$ws.Rows.Item(179).Insert()

# Fill in the new record in row 179.
$ws.Range("A179").Value = 10
$ws.Range("B179").Value = "Vega Modelo de Temuco"
$ws.Range("C179").Value = "La Araucanía"
$ws.Range("D179").Value = 45009
$ws.Range("D179").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E179").Value = 9
$ws.Range("F179").Value = "Fruta"
$ws.Range("G179").Value = 100104
$ws.Range("H179").Value = "Frutos de pepita"
$ws.Range("I179").Value = 100104003
$ws.Range("J179").Value = "Membrillo"
$ws.Range("K179").Value = "Champion"
$ws.Range("L179").Value = "Primera"
$ws.Range("M179").Value = 80
$ws.Range("N179").Value = 14000
$ws.Range("O179").Value = 14000
$ws.Range("P179").Value = 14000
$ws.Range("Q179").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R179").Value = "Región de O'Higgins"
$ws.Range("S179").Value = 778
$ws.Range("T179").Value = 18
